$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 is a brand-new timesheet entry for the same day as row 11.
# Clone the date/time cell formatting from row 11 (B11/C11) before row 11's
# own contents are edited, then fill in the actual values for row 12.
$ws.Range("B11").Copy($ws.Range("B12"))
$ws.Range("B12").Value2 = 43503

$ws.Range("C11").Copy($ws.Range("C12"))
$ws.Range("C12").Value2 = 0.85416666666666663

# Row 11: A11 now records 1.5 hours worked, and C11 switches from a bare time
# value to the descriptive text "16:30 - 18:00".
$ws.Range("A11").Value = 1.5
$ws.Range("C11").Value = "16:30 - 18:00"

# The active selection moved to C12.
$ws.Range("C12").Select()
